# update of GL for face images
# The CNN sheet's trial table grows from 10 trials (rows 3:12) to 20 trials
# (rows 3:22); the MEAN/STD summary formulas move down from rows 13:14 to
# rows 23:24; the header labels shift from standalone rows 17:18 up into
# C1:C2; and the "window size" helper row (D7:J9) gets new values. The two
# charts on GL_adaptive then need their category range to track the last
# five "window size" columns (F:J instead of D:H) and their GL_MV /
# GL_adaptive value ranges shift from column B to column D.

$wb = $excel.ActiveWorkbook

$cnn = $wb.Worksheets.Item("CNN")
$glmv = $wb.Worksheets.Item("GL_MV")
$gladaptive = $wb.Worksheets.Item("GL_adaptive")

# ---- Headers -------------------------------------------------------------
$cnn.Range("A1").Value = "CNN"
$cnn.Range("C1").Value = "Task: L vs D"
$cnn.Range("A2").Value = "training_acc"
$cnn.Range("B2").Value = "Test_acc"
$cnn.Range("C2").Value = "15 Letter (J&E) and 15 digit (5&7) in training"

# ---- Per-trial data (20 trials instead of 10) -----------------------------
# (this also overwrites the old standalone text labels that used to live at
# A17/A18 - they move into C1/C2 above.)
$AValues = @(
    0,
    0,
    0,
    0,
    0,
    0,
    0.033333300000000003,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)
$BValues = @(
    0.28199999999999997,
    0.28100000000000003,
    0.28699999999999998,
    0.27700000000000002,
    0.246,
    0.27300000000000002,
    0.28699999999999998,
    0.29399999999999998,
    0.29299999999999998,
    0.30399999999999999,
    0.28100000000000003,
    0.29099999999999998,
    0.28599999999999998,
    0.29899999999999999,
    0.28899999999999998,
    0.28399999999999997,
    0.26800000000000002,
    0.26700000000000002,
    0.28199999999999997,
    0.27600000000000002
)

for ($i = 0; $i -lt $AValues.Length; $i++) {
    $row = 3 + $i
    $cnn.Cells.Item($row, 1).Value = $AValues[$i]
    $cnn.Cells.Item($row, 2).Value = $BValues[$i]
}

# The old "mean"/"STD" row labels (C13/C14) are no longer needed there -
# they get re-created further down at C23/C24.
$cnn.Range("C13:C14").ClearContents()

# ---- Window-size helper block (D7:J9) -------------------------------------
$windowSizes = @(5, 10, 15, 18, 21, 24, 27)
for ($i = 0; $i -lt $windowSizes.Length; $i++) {
    $col = 4 + $i
    $cnn.Cells.Item(7, $col).Value = $windowSizes[$i]
    $cnn.Cells.Item(8, $col).Value = 0.28234999999999999
    $cnn.Cells.Item(9, $col).Value = 0.0016666650000000001
}

# ---- Summary formulas now live at rows 23:24 instead of 13:14 ------------
$cnn.Range("A23").Formula = "=AVERAGE(A3:A22)"
$cnn.Range("B23").Formula = "=AVERAGE(B3:B22)"
$cnn.Range("C23").Value = "mean"

$cnn.Range("A24").Formula = "=STDEV(A3:A22)"
$cnn.Range("B24").Formula = "=STDEV(B3:B22)"
$cnn.Range("C24").Value = "STD"

# ---- Chart source ranges on GL_adaptive -----------------------------------
# Chart 1 ("Test error"): categories move from D7:J7 to F7:J7, and the
# GL_MV/GL_adaptive series move from column B to column D.
$chart1 = $gladaptive.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = '=SERIES("CNN+GL(MV)",CNN!$F$7:$J$7,GL_MV!$D$28:$H$28,1)'
$chart1.SeriesCollection().Item(2).Formula = '=SERIES("CNN+GL(Adaptive)",CNN!$F$7:$J$7,GL_adaptive!$D$28:$H$28,2)'
$chart1.SeriesCollection().Item(3).Formula = '=SERIES("CNN",''[1]GL_adaptive''!$A$11:$J$11,CNN!$F$8:$J$8,3)'

# Chart 2 ("Training error"): categories move from D7:H7 to F7:J7, and the
# GL_MV/GL_adaptive series move from column B to column D.
$chart2 = $gladaptive.ChartObjects().Item(2).Chart
$chart2.SeriesCollection().Item(1).Formula = '=SERIES("CNN+GL(MV)",CNN!$F$7:$J$7,GL_MV!$D$13:$H$13,1)'
$chart2.SeriesCollection().Item(2).Formula = '=SERIES("CNN+GL(Adaptive)",CNN!$F$7:$J$7,GL_adaptive!$D$13:$H$13,2)'
$chart2.SeriesCollection().Item(3).Formula = '=SERIES("CNN",''[1]GL_adaptive''!$A$11:$J$11,CNN!$F$9:$J$9,3)'
